$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 15.51809326405486
$ws.Range("C2").Value = 11.44037335774968
$ws.Range("D2").Value = 5.136884144943457
$ws.Range("F2").Value = 27.18648119850263
$ws.Range("G2").Value = 34.07047412151667
$ws.Range("H2").Value = 15.35671303431407
$ws.Range("I2").Value = 23.87802735514514
$ws.Range("L2").Value = 10.88501830446585
$ws.Range("M2").Value = 15.25745719148344
$ws.Range("B3").Value = 14.97142024420893
$ws.Range("C3").Value = 11.07997703539071
$ws.Range("D3").Value = 5.135804338773383
$ws.Range("F3").Value = 27.14869940853782
$ws.Range("G3").Value = 33.96289227374837
$ws.Range("H3").Value = 15.4024602874613
$ws.Range("I3").Value = 23.99662326291098
$ws.Range("L3").Value = 10.90430630106511
$ws.Range("M3").Value = 15.1571382546325
$ws.Range("B4").Value = 14.62759853051488
$ws.Range("C4").Value = 10.8507329561723
$ws.Range("D4").Value = 5.135383349924752
$ws.Range("F4").Value = 27.13541968662009
$ws.Range("G4").Value = 33.91198337104156
$ws.Range("H4").Value = 15.43436440037365
$ws.Range("I4").Value = 24.07591398217671
$ws.Range("L4").Value = 10.91780819101911
$ws.Range("M4").Value = 15.09772482070409
$ws.Range("B5").Value = 14.48565241196841
$ws.Range("C5").Value = 10.75539930983479
$ws.Range("D5").Value = 5.135273086780199
$ws.Range("F5").Value = 27.1325012903671
$ws.Range("G5").Value = 33.89505348473875
$ws.Range("H5").Value = 15.44832146178543
$ws.Range("I5").Value = 24.10984664296998
$ws.Range("L5").Value = 10.9237274170063
$ws.Range("M5").Value = 15.07408172517825
$ws.Range("B6").Value = 14.46197820127997
$ws.Range("C6").Value = 10.73945627618256
$ws.Range("D6").Value = 5.135258492401841
$ws.Range("F6").Value = 27.13216719284789
$ws.Range("G6").Value = 33.89247284024806
$ws.Range("H6").Value = 15.45069664946511
$ws.Range("I6").Value = 24.11557885349368
$ws.Range("L6").Value = 10.92473548870501
$ws.Range("M6").Value = 15.07019068698651
$ws.Range("B7").Value = 14.62569132843229
$ws.Range("C7").Value = 10.84945488568366
$ws.Range("D7").Value = 5.135381614142204
$ws.Range("F7").Value = 27.13537023665776
$ws.Range("G7").Value = 33.91173959380664
$ws.Range("H7").Value = 15.4345487646158
$ws.Range("I7").Value = 24.07636505496617
$ws.Range("L7").Value = 10.91788633108775
$ws.Range("M7").Value = 15.09740363530315
$ws.Range("B8").Value = 15.33141813020069
$ws.Range("C8").Value = 11.31781001867133
$ws.Range("D8").Value = 5.136461857555289
$ws.Range("F8").Value = 27.17139581816425
$ws.Range("G8").Value = 34.03024246483425
$ws.Range("H8").Value = 15.37169282183989
$ws.Range("I8").Value = 23.91757199616979
$ws.Range("L8").Value = 10.89132445469653
$ws.Range("M8").Value = 15.22242612164275
$ws.Range("B9").Value = 16.6420903519607
$ws.Range("C9").Value = 12.16971695674124
$ws.Range("D9").Value = 5.140481885699981
$ws.Range("F9").Value = 27.32064410130374
$ws.Range("G9").Value = 34.38217952831447
$ws.Range("H9").Value = 15.27885253832042
$ws.Range("I9").Value = 23.65782371458257
$ws.Range("L9").Value = 10.85240212260429
$ws.Range("M9").Value = 15.48402163854168
$ws.Range("B10").Value = 17.55039123847094
$ws.Range("C10").Value = 12.75108079583391
$ws.Range("D10").Value = 5.144570498975904
$ws.Range("F10").Value = 27.47788913574175
$ws.Range("G10").Value = 34.71233559569723
$ws.Range("H10").Value = 15.22938179462546
$ws.Range("I10").Value = 23.49886566058971
$ws.Range("L10").Value = 10.83183341438739
$ws.Range("M10").Value = 15.68498756348564
$ws.Range("B11").Value = 17.94999422771518
$ws.Range("C11").Value = 13.00519650093443
$ws.Range("D11").Value = 5.146671894641316
$ws.Range("F11").Value = 27.55962550370431
$ws.Range("G11").Value = 34.87766654007971
$ws.Range("H11").Value = 15.21098314139789
$ws.Range("I11").Value = 23.43355845307427
$ws.Range("L11").Value = 10.82421868298168
$ws.Range("M11").Value = 15.77804046396367
$ws.Range("B12").Value = 18.09923736235928
$ws.Range("C12").Value = 13.09988620391706
$ws.Range("D12").Value = 5.147501920223506
$ws.Range("F12").Value = 27.59202847229164
$ws.Range("G12").Value = 34.94240674668719
$ws.Range("H12").Value = 15.20460924093188
$ws.Range("I12").Value = 23.40984245747587
$ws.Range("L12").Value = 10.82158553677403
$ws.Range("M12").Value = 15.81348800425394
$ws.Range("B13").Value = 18.06718940588647
$ws.Range("C13").Value = 13.07956223697213
$ws.Range("D13").Value = 5.147321642230218
$ws.Range("F13").Value = 27.58498566136155
$ws.Range("G13").Value = 34.9283696860259
$ws.Range("H13").Value = 15.20595554545515
$ws.Range("I13").Value = 23.41490488838438
$ws.Range("L13").Value = 10.82214149842102
$ws.Range("M13").Value = 15.805844773497
$ws.Range("B14").Value = 17.96231484850924
$ws.Range("C14").Value = 13.01301777724835
$ws.Range("D14").Value = 5.146739497045792
$ws.Range("F14").Value = 27.56226233508369
$ws.Range("G14").Value = 34.88295028191004
$ws.Range("H14").Value = 15.2104468511804
$ws.Range("I14").Value = 23.43158695142988
$ws.Range("L14").Value = 10.82399703522153
$ws.Range("M14").Value = 15.78095267743832
$ws.Range("B15").Value = 17.89780214249605
$ws.Range("C15").Value = 12.97205564495358
$ws.Range("D15").Value = 5.146387366280524
$ws.Range("F15").Value = 27.54853209586108
$ws.Range("G15").Value = 34.85540594873023
$ws.Range("H15").Value = 15.21327524790968
$ws.Range("I15").Value = 23.44193752162397
$ws.Range("L15").Value = 10.82516620833441
$ws.Range("M15").Value = 15.76573222711166
$ws.Range("B16").Value = 17.52399187177466
$ws.Range("C16").Value = 12.7342610500715
$ws.Range("D16").Value = 5.144437988265721
$ws.Range("F16").Value = 27.47275149562188
$ws.Range("G16").Value = 34.70183185128921
$ws.Range("H16").Value = 15.23066705893613
$ws.Range("I16").Value = 23.50327519674962
$ws.Range("L16").Value = 10.83236610097536
$ws.Range("M16").Value = 15.67893711472047
$ws.Range("B17").Value = 17.29109805084362
$ws.Range("C17").Value = 12.58569334423342
$ws.Range("D17").Value = 5.14330363287964
$ws.Range("F17").Value = 27.42886561243207
$ws.Range("G17").Value = 34.61146744200281
$ws.Range("H17").Value = 15.24239007605806
$ws.Range("I17").Value = 23.54270288015629
$ws.Range("L17").Value = 10.83722912143433
$ws.Range("M17").Value = 15.62609203294579
$ws.Range("B18").Value = 17.15587112077457
$ws.Range("C18").Value = 12.49926956313348
$ws.Range("D18").Value = 5.142673936304083
$ws.Range("F18").Value = 27.40458530047623
$ws.Range("G18").Value = 34.56092042238721
$ws.Range("H18").Value = 15.24951914879196
$ws.Range("I18").Value = 23.56603941736979
$ws.Range("L18").Value = 10.8401901855508
$ws.Range("M18").Value = 15.59585207481841
$ws.Range("B19").Value = 17.10987107614164
$ws.Range("C19").Value = 12.46984271701353
$ws.Range("D19").Value = 5.142464653720694
$ws.Range("F19").Value = 27.39653002227241
$ws.Range("G19").Value = 34.54405262187545
$ws.Range("H19").Value = 15.25199918775395
$ws.Range("I19").Value = 23.57405369010121
$ws.Range("L19").Value = 10.84122091804974
$ws.Range("M19").Value = 15.58564072560956
$ws.Range("B20").Value = 17.3160226486061
$ws.Range("C20").Value = 12.6016095661902
$ws.Range("D20").Value = 5.143422035170856
$ws.Range("F20").Value = 27.43343791703345
$ws.Range("G20").Value = 34.62093937749754
$ws.Range("H20").Value = 15.24110213850008
$ws.Range("I20").Value = 23.5384375018985
$ws.Range("L20").Value = 10.83669447363758
$ws.Range("M20").Value = 15.6317016015011
$ws.Range("B21").Value = 17.99317635171287
$ws.Range("C21").Value = 13.03260561938277
$ws.Range("D21").Value = 5.146909560628739
$ws.Range("F21").Value = 27.5688974862014
$ws.Range("G21").Value = 34.89623355492459
$ws.Range("H21").Value = 15.20911152296616
$ws.Range("I21").Value = 23.42665943573012
$ws.Range("L21").Value = 10.82344522489499
$ws.Range("M21").Value = 15.78825857368417
$ws.Range("B22").Value = 18.42357728854222
$ws.Range("C22").Value = 13.30530040143541
$ws.Range("D22").Value = 5.149388419036581
$ws.Range("F22").Value = 27.66587716345277
$ws.Range("G22").Value = 35.08856233624428
$ws.Range("H22").Value = 15.19166300634274
$ws.Range("I22").Value = 23.35952223301084
$ws.Range("L22").Value = 10.81624542688689
$ws.Range("M22").Value = 15.89179268635349
$ws.Range("B23").Value = 18.19501370450577
$ws.Range("C23").Value = 13.16059530907914
$ws.Range("D23").Value = 5.148047294483363
$ws.Range("F23").Value = 27.61335031574313
$ws.Range("G23").Value = 34.9847935784234
$ws.Range("H23").Value = 15.2006582050225
$ws.Range("I23").Value = 23.39481087786626
$ws.Range("L23").Value = 10.81995461869264
$ws.Range("M23").Value = 15.8364315003934
$ws.Range("B24").Value = 17.30475838091337
$ws.Range("C24").Value = 12.59441698512873
$ws.Range("D24").Value = 5.143368435473232
$ws.Range("F24").Value = 27.43136781681746
$ws.Range("G24").Value = 34.61665273657603
$ws.Range("H24").Value = 15.24168320240631
$ws.Range("I24").Value = 23.54036379585026
$ws.Range("L24").Value = 10.83693567316108
$ws.Range("M24").Value = 15.62916507329436
$ws.Range("B25").Value = 16.29646668724314
$ws.Range("C25").Value = 11.94681583648768
$ws.Range("D25").Value = 5.13919327603242
$ws.Range("F25").Value = 27.27187519078371
$ws.Range("G25").Value = 34.27428032389565
$ws.Range("H25").Value = 15.30069085385428
$ws.Range("I25").Value = 23.72252279455027
$ws.Range("L25").Value = 10.86152183607969
$ws.Range("M25").Value = 15.41162029888642
